# vendas.xlsx re-upload:
#  - worksheet renamed "Sheet 1" -> "minhas_vendas"
#  - table "Tabela2" renamed (Name + DisplayName) -> "vendas"
#  - stray formatted-but-empty cell H84 removed (it carried no value, only the
#    now-unused underline style), which also shrinks the sheet dimension from
#    A1:H88 back down to A1:F88
#  - the sheet's saved selection now covers the data range A2:F88

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "minhas_vendas"

# Rename the workbook's table (both the internal Name and the DisplayName
# track together through the ListObject.Name setter).
$tbl = $ws.ListObjects.Item(1)
$tbl.Name = "vendas"

# H84 held no value, just a leftover cell format -- clear it out entirely so
# it disappears from the sheet (not merely blanked), which also pulls the
# sheet's used range/dimension back in from H88 to F88.
$ws.Range("H84").Clear()

# Leave the sheet with the data body selected, as last saved.
$ws.Range("A2:F88").Select()
